# Apply updated "想去人数" (F column) figures to the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Row -> new value for sheet "展览"
$updates1 = @{
    5  = 3332
    6  = 1093
    7  = 2224
    12 = 1685
    17 = 227
    18 = 1599
    19 = 643
    20 = 734
    22 = 12296
    23 = 12356
    25 = 706
    29 = 382
    30 = 1928
}

foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# Row -> new value for sheet "全部类型"
$updates4 = @{
    6  = 3332
    7  = 1093
    8  = 2224
    13 = 1685
    21 = 227
    22 = 1599
    23 = 643
    24 = 734
    26 = 12296
    27 = 12356
    29 = 706
    33 = 382
    34 = 1928
}

foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
